$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 563.3333
$ws.Range("J41").Value = 1490
$ws.Range("L41").Value = 1490
$ws.Range("N41").Value = -2370

$ws.Range("H116").Value = 1811.8823
$ws.Range("I116").Value = 1628.76
$ws.Range("J116").Value = 2320.5557
$ws.Range("K116").Value = 1628.76
$ws.Range("L116").Value = 2320.5557
$ws.Range("M116").Value = 1813.24
$ws.Range("N116").Value = -9204.555700000001

$ws.Range("H137").Value = 1521.9403
$ws.Range("I137").Value = 1132.4681
$ws.Range("J137").Value = 2437.2
$ws.Range("K137").Value = 3397.4043
$ws.Range("L137").Value = 7311.599999999999
$ws.Range("M137").Value = -847.4043000000001
$ws.Range("N137").Value = -12411.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1369.8
$ws.Range("I45").Value = 1294.1111
$ws.Range("J45").Value = 1710.4
$ws.Range("K45").Value = 1294.1111
$ws.Range("L45").Value = 1710.4
$ws.Range("M45").Value = -917.1111000000001
$ws.Range("N45").Value = -2464.4

$ws.Range("H61").Value = 6505.6035
$ws.Range("I61").Value = 3510.7234
$ws.Range("J61").Value = 19301.908
$ws.Range("K61").Value = 3510.7234
$ws.Range("L61").Value = 19301.908
$ws.Range("M61").Value = -3298.7234
$ws.Range("N61").Value = -19725.908

$ws.Range("H74").Value = 3933.2927
$ws.Range("I74").Value = 1701.9354
$ws.Range("J74").Value = 10850.5
$ws.Range("K74").Value = 1701.9354
$ws.Range("L74").Value = 10850.5
$ws.Range("M74").Value = -827.9354000000001
$ws.Range("N74").Value = -12598.5

$ws.Range("H77").Value = 3933.2927
$ws.Range("I77").Value = 1701.9354
$ws.Range("J77").Value = 10850.5
$ws.Range("K77").Value = 8509.677
$ws.Range("L77").Value = 54252.5
$ws.Range("M77").Value = -4141.677
$ws.Range("N77").Value = -62988.5

$ws.Range("H122").Value = 1691.9565
$ws.Range("I122").Value = 1309.2
$ws.Range("J122").Value = 2409.625
$ws.Range("K122").Value = 3927.6
$ws.Range("L122").Value = 7228.875
$ws.Range("M122").Value = -1477.6
$ws.Range("N122").Value = -12128.875

$ws.Range("H132").Value = 2233
$ws.Range("I132").Value = 1901.72
$ws.Range("J132").Value = 4303.5
$ws.Range("K132").Value = 5705.16
$ws.Range("L132").Value = 12910.5
$ws.Range("M132").Value = -3175.16
$ws.Range("N132").Value = -17970.5

$ws.Range("H136").Value = 6505.6035
$ws.Range("I136").Value = 3510.7234
$ws.Range("J136").Value = 19301.908
$ws.Range("K136").Value = 10532.1702
$ws.Range("L136").Value = 57905.724
$ws.Range("M136").Value = -7982.1702
$ws.Range("N136").Value = -63005.724

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H55").Value = 60900
$ws.Range("J55").Value = 60900
$ws.Range("L55").Value = 60900
$ws.Range("N55").Value = -61446

$ws.Range("H105").Value = 981296.6
$ws.Range("I105").Value = 1363194.9
$ws.Range("J105").Value = 5334.4443
$ws.Range("K105").Value = 1363194.9
$ws.Range("L105").Value = 5334.4443
$ws.Range("M105").Value = -1361447.9
$ws.Range("N105").Value = -8828.444299999999

$ws.Range("H107").Value = 1540
$ws.Range("I107").Value = 1350
$ws.Range("K107").Value = 1350
$ws.Range("M107").Value = 570

$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()

$ws.Range("H134").Value = 17050.422
$ws.Range("I134").Value = 1168.3529
$ws.Range("J134").Value = 79357
$ws.Range("K134").Value = 3505.0587
$ws.Range("L134").Value = 238071
$ws.Range("M134").Value = -970.0587000000005
$ws.Range("N134").Value = -243141

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 160
$ws.Range("I22").Value = 153.75
$ws.Range("K22").Value = 153.75
$ws.Range("M22").Value = 196.25

$ws.Range("H31").Value = 1985.6608
$ws.Range("I31").Value = 1362.1086
$ws.Range("K31").Value = 1362.1086
$ws.Range("M31").Value = -1067.1086

$ws.Range("H34").Value = 1985.6608
$ws.Range("I34").Value = 1362.1086
$ws.Range("K34").Value = 1362.1086
$ws.Range("M34").Value = -1160.1086

$ws.Range("H99").Value = 3955.3635
$ws.Range("I99").Value = 2861.8
$ws.Range("J99").Value = 4866.6665
$ws.Range("K99").Value = 2861.8
$ws.Range("L99").Value = 4866.6665
$ws.Range("M99").Value = -1363.8
$ws.Range("N99").Value = -7862.6665

$ws.Range("H105").Value = 1348.0667
$ws.Range("I105").Value = 828.5714
$ws.Range("J105").Value = 1802.625
$ws.Range("K105").Value = 828.5714
$ws.Range("L105").Value = 1802.625
$ws.Range("M105").Value = 918.4286
$ws.Range("N105").Value = -5296.625

$ws.Range("H107").Value = 1275.6316
$ws.Range("I107").Value = 1427.8334
$ws.Range("J107").Value = 1014.7143
$ws.Range("K107").Value = 1427.8334
$ws.Range("L107").Value = 1014.7143
$ws.Range("M107").Value = 492.1666
$ws.Range("N107").Value = -4854.7143

$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()

$ws.Range("H126").Value = 3955.3635
$ws.Range("I126").Value = 2861.8
$ws.Range("J126").Value = 4866.6665
$ws.Range("K126").Value = 8585.400000000001
$ws.Range("L126").Value = 14599.9995
$ws.Range("M126").Value = -6115.400000000001
$ws.Range("N126").Value = -19539.9995

$ws.Range("H132").Value = 2681
$ws.Range("I132").Value = 2842.9678
$ws.Range("J132").Value = 2224.5454
$ws.Range("K132").Value = 8528.903399999999
$ws.Range("L132").Value = 6673.6362
$ws.Range("M132").Value = -5998.903399999999
$ws.Range("N132").Value = -11733.6362

$ws.Range("H134").Value = 2211.3167
$ws.Range("I134").Value = 1239.7354
$ws.Range("J134").Value = 3481.8462
$ws.Range("K134").Value = 3719.2062
$ws.Range("L134").Value = 10445.5386
$ws.Range("M134").Value = -1184.2062
$ws.Range("N134").Value = -15515.5386

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 473.27274
$ws.Range("J92").Value = 487.33334
$ws.Range("L92").Value = 1462.00002
$ws.Range("N92").Value = -3958.00002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6234.162
$ws.Range("I70").Value = 5897.231
$ws.Range("J70").Value = 6416.6665
$ws.Range("K70").Value = 5897.231
$ws.Range("L70").Value = 6416.6665
$ws.Range("M70").Value = -5627.231
$ws.Range("N70").Value = -6956.6665

$ws.Range("H73").Value = 6234.162
$ws.Range("I73").Value = 5897.231
$ws.Range("J73").Value = 6416.6665
$ws.Range("K73").Value = 5897.231
$ws.Range("L73").Value = 6416.6665
$ws.Range("M73").Value = -4961.231
$ws.Range("N73").Value = -8288.666499999999

$ws.Range("H80").Value = 3177
$ws.Range("I80").Value = 2005
$ws.Range("J80").Value = 3274.6667
$ws.Range("K80").Value = 2005
$ws.Range("L80").Value = 3274.6667
$ws.Range("M80").Value = -1007
$ws.Range("N80").Value = -5270.6667

$ws.Range("H83").Value = 3177
$ws.Range("I83").Value = 2005
$ws.Range("J83").Value = 3274.6667
$ws.Range("K83").Value = 10025
$ws.Range("L83").Value = 16373.3335
$ws.Range("M83").Value = -5033
$ws.Range("N83").Value = -26357.3335

$ws.Range("H102").Value = 3396.2632
$ws.Range("I102").Value = 3236.72
$ws.Range("J102").Value = 3703.077
$ws.Range("K102").Value = 3236.72
$ws.Range("L102").Value = 3703.077
$ws.Range("M102").Value = -1614.72
$ws.Range("N102").Value = -6947.077

$ws.Range("H107").Value = 220.125
$ws.Range("I107").Value = 135.57143
$ws.Range("J107").Value = 285.8889
$ws.Range("K107").Value = 135.57143
$ws.Range("L107").Value = 285.8889
$ws.Range("M107").Value = 1784.42857
$ws.Range("N107").Value = -4125.8889

$ws.Range("H112").Value = 79800
$ws.Range("J112").Value = 79800
$ws.Range("L112").Value = 79800
$ws.Range("N112").Value = -82016

$ws.Range("H132").Value = 6157.3335
$ws.Range("I132").Value = 4121.892
$ws.Range("J132").Value = 21219.6
$ws.Range("K132").Value = 12365.676
$ws.Range("L132").Value = 63658.8
$ws.Range("M132").Value = -9835.675999999999
$ws.Range("N132").Value = -68718.79999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1092.4615
$ws.Range("I22").Value = 963.6667
$ws.Range("J22").Value = 1202.8572
$ws.Range("K22").Value = 963.6667
$ws.Range("L22").Value = 1202.8572
$ws.Range("M22").Value = -668.6667
$ws.Range("N22").Value = -1792.8572

$ws.Range("H27").Value = 1092.4615
$ws.Range("I27").Value = 963.6667
$ws.Range("J27").Value = 1202.8572
$ws.Range("K27").Value = 963.6667
$ws.Range("L27").Value = 1202.8572
$ws.Range("M27").Value = -856.6667
$ws.Range("N27").Value = -1416.8572

$ws.Range("H46").Value = 1104.5714
$ws.Range("I46").Value = 975
$ws.Range("J46").Value = 1156.4
$ws.Range("K46").Value = 975
$ws.Range("L46").Value = 1156.4
$ws.Range("M46").Value = -787
$ws.Range("N46").Value = -1532.4

$ws.Range("H94").Value = 24000
$ws.Range("J94").Value = 24000
$ws.Range("L94").Value = 24000
$ws.Range("N94").Value = -25352

$ws.Range("H112").Value = 79800
$ws.Range("J112").Value = 79800
$ws.Range("L112").Value = 79800
$ws.Range("N112").Value = -82754

$ws.Range("H122").Value = 6646.3066
$ws.Range("I122").Value = 6276.095
$ws.Range("J122").Value = 7423.75
$ws.Range("K122").Value = 18828.285
$ws.Range("L122").Value = 22271.25
$ws.Range("M122").Value = -16378.285
$ws.Range("N122").Value = -27171.25

$ws.Range("H132").Value = 2613.8408
$ws.Range("I132").Value = 2250.725
$ws.Range("K132").Value = 6752.174999999999
$ws.Range("M132").Value = -4222.174999999999

$ws.Range("H136").Value = 2727.553
$ws.Range("I136").Value = 1507.5344
$ws.Range("J136").Value = 5348.3335
$ws.Range("K136").Value = 4522.6032
$ws.Range("L136").Value = 16045.0005
$ws.Range("M136").Value = -1972.6032
$ws.Range("N136").Value = -21145.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 816.9452
$ws.Range("I132").Value = 258.23727
$ws.Range("J132").Value = 3171.5
$ws.Range("K132").Value = 774.71181
$ws.Range("L132").Value = 9514.5
$ws.Range("M132").Value = 1755.28819
$ws.Range("N132").Value = -14574.5
